# Graphics II Project Rubric.xlsx edit
# "Put three optional shaders runtime compiling function in to the D3DUtils"
#
# Marks three more rubric rows as completed on Milestone III (column E = "III",
# column F = "X"), and marks the Milestone III GIT/API-cleanup carry-over rows
# (90/91, columns D & E) with "X" as well. Also updates the active view's
# selection/scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 8  - "All Graphics API Objects cleaned up in memory" -> Milestone III, complete
$ws.Range("E8").Value = "III"
$ws.Range("F8").Value = "X"

# Row 19 - Geometry Instancing w/ 1 drawInstance related call -> Milestone III, complete
$ws.Range("E19").Value = "III"
$ws.Range("F19").Value = "X"

# Row 21 - Advanced use textures or mapping : Shadow Mapping -> Milestone III, complete
$ws.Range("E21").Value = "III"
$ws.Range("F21").Value = "X"

# Row 38 - Advanced Lighting : Shadow Mapping -> Milestone III, complete
$ws.Range("E38").Value = "III"
$ws.Range("F38").Value = "X"

# Row 47 - Substantial Use of Compute Shader (object transformations) -> Milestone III, complete
$ws.Range("E47").Value = "III"
$ws.Range("F47").Value = "X"

# Rows 90 & 91 - Effective use of GIT / cleaned up API objects now also
# carried over onto Milestone III (column E), in addition to Milestone II
# (column D).
$ws.Range("D90").Value = "X"
$ws.Range("E90").Value = "X"
$ws.Range("D91").Value = "X"
$ws.Range("E91").Value = "X"

# Update the sheet view: scroll so row 13 / column B is the top-left visible
# cell, and select F20.
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 2
$ws.Range("F20").Select()
